# payments.xlsx: add payment 79174445 (Cash) 2025-08-29T16:20:43

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6, column A ("phone") was stored as text; normalize it to a real number,
# matching the other rows in the sheet.
$ws.Range("A6").Value = 79174445

# Append the new payment as row 7.
# Column A ("phone") must stay textual even though it looks numeric, so
# prefix with an apostrophe to force text entry, then clear the resulting
# "quote prefix" cell style so no extra formatting is left behind.
$ws.Range("A7").Value = "'79174445"
$ws.Range("A7").Style = "Normal"

$ws.Range("B7").Value = 5
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 5
$ws.Range("G7").Value = "Cash"
$ws.Range("H7").Value = "2025-08-29T16:20:43"
